# Generate Report for Handoff
# Update the "423069c5-47e9-41c5-a0d5-9576d90d397d.md" row to reflect that the
# file is now ready for handoff (status changed, new handoff timestamps, and
# a version-mismatch error detail recorded for each locale).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fd4f8e33c8c0e657ac21b1cec07c862c304d125/e2e/423069c5-47e9-41c5-a0d5-9576d90d397d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d84232566f983770e06ed58c4be92f961a8d196/e2e/423069c5-47e9-41c5-a0d5-9576d90d397d.md."

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 18:55:28"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-27 18:55:23"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-27 18:55:28"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
